$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-blank "Outcome" cells for already-existing rows ---
$ws.Range("D11").Value = "Ongoing"
$ws.Range("D12").Value = "Ongoing"
$ws.Range("D13").Value = "Ongoing"
$ws.Range("D14").Value = "Ongoing"

# --- Row 15: Poloniex, LLC settlement ---
$ws.Range("D15").Value = "Settlment"
$ws.Range("E15").Value = "Unregistered Exchange"
$ws.Range("M15").Value = "Washington, D.C."
$ws.Range("H15").Value = "Poloniex"
$ws.Range("F15").Value = "Civil"
$ws.Range("G15").Value = "N/A"
$ws.Range("I15").Value = "N/A"
$ws.Range("J15").Value = 10000000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1

# --- Row 16: Blockchain Credit Partners d/b/a DeFi Money Market ---
$ws.Range("H16").Value = "DeFi Money Market"
$ws.Range("G16").Value = "mTokens"
$ws.Range("D16").Value = "Settlement"
$ws.Range("E16").Value = "Unregistered Offering"
$ws.Range("F16").Value = "Civil"
$ws.Range("I16").Value = "Ethereum"
$ws.Range("M16").Value = "New York"
$ws.Range("J16").Value = 30000000
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1

# --- Row 17: SEC v. Uulala, Inc., et al. ---
$ws.Range("G17").Value = "UULA"
$ws.Range("H17").Value = "Uulala, Inc."
$ws.Range("D17").Value = "Ongoing"
$ws.Range("E17").Value = "Unregistered Offering"
$ws.Range("F17").Value = "Civil"
$ws.Range("I17").Value = "N/A"
$ws.Range("M17").Value = "Los Angeles"
$ws.Range("J17").Value = 9000000
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 1

# --- Row 18: Blotics Ltd., f/d/b/a Coinschedule Ltd. ---
$ws.Range("E18").Value = "Failed Disclosure"
$ws.Range("H18").Value = "Coinschedule"
$ws.Range("D18").Value = "Settlement"
$ws.Range("F18").Value = "Civil"
$ws.Range("G18").Value = "N/A"
$ws.Range("I18").Value = "N/A"
$ws.Range("M18").Value = "Washington, D.C."
$ws.Range("J18").Value = 197000
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("N18").Select()
